$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old hidden "_GoBack" bookmark that currently sits
#    between "officer" and ".name}}," in the greeting line.
#    (It will be re-added further down, next to the new stats text.)
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the "[insert stats here]" placeholder paragraph and replace
#    it (and the content that follows) with the real stats section:
#       STATS:
#       Since going live, ProspectAve.io has accumulated over 6,800
#         pageviews and 2,300 unique users   <-- bookmark goes here
#       Active monthly users peaked at 1,223 over LWNPARTIES weekend
#       25% of our users access the page through their phones
#
#    Using InsertXML lets us place the bookmark start/end tags exactly
#    where they belong (immediately after "...unique users", inside
#    that paragraph) instead of relying on Range/Bookmarks.Add with a
#    collapsed range sitting on a paragraph boundary.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "[insert stats here]") {
        $target = $para
        break
    }
}

$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>STATS:</w:t></w:r></w:p>' + `
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Since going live, ProspectAve.io has accumulated over 6,800 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pageviews</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and 2,300 unique users</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Active</w:t></w:r><w:r><w:t xml:space="preserve"> monthly users peaked at 1,223 over LWNPARTIES weekend</w:t></w:r></w:p>' + `
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>25% of our users access the page through their phones</w:t></w:r></w:p>'

$target.Range.InsertXML($xmlFragment)
